# Review_119.docx -> "Review 118" content swap
#
# Strategy: use Find.Execute purely to *locate* each distinct run of text
# (replacement text = "" and Wrap = wdFindStop / 0, i.e. no in-place
# replace performed by Find itself), then assign the new text straight to
# the found Range's .Text property. Setting .Text directly (instead of
# passing the replacement through Find.Execute's Replace argument) avoids
# Word's AutoCorrect/AutoFormat "smart quotes" pass mangling the Hebrew
# text's straight apostrophes ('). A manual line break (a fresh <w:br/>)
# is inserted into replacement text with [char]11 (vertical tab, the same
# control character wdFindContinue/Range.Text uses for ^l).

$d = $word.ActiveDocument
$lb = [char]11   # -> <w:br/>

function Set-RangeText($old, $new) {
    $rng = $d.Content
    $found = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Output "NOT FOUND: $old"
        return
    }
    $rng.Text = $new
}

# --- Title heading block (Heading1 paragraph, single run, 2 <w:t> split by <w:br/>) ---
Set-RangeText "Review 119: [Short] FLIRT: Feedback Loop In-context Red Teaming, 09.08.23" "Review 118: [Short] Seeing through the Brain: Image Reconstruction of Visual Perception from Human Brain Signals, 08.08.23"

Set-RangeText "https://huggingface.co/papers/2308.04265" "https://huggingface.co/papers/2308.02510"

# --- Bold "Paper:" line ------------------------------------------------------
Set-RangeText "Paper: https://arxiv.org/abs/2308.04265v2" "Paper: https://arxiv.org/abs/2208.03666v4"

# --- Body paragraph (single run, four <w:t> separated by <w:br/><w:br/>) ----

Set-RangeText "בטיחות מודלים גנרטיביים הינו אחד מנושאי המחקר החמים בבינה מלאכותית גנרטיבית (GenAI). הרי אנחנו לא רוצים מודל המצייר תמונה לפי התיאור הטקסטואלי יגנרט לנו תמונה קשה, אלימה או מטרידה גם אם נבקש את זה ממנו. למניעת תופעות אלו צריך לזהות פרומפטים מתוחכמים שגורמים למודל ליצור תוכן בעייתי. " "מכונה שיודעת לקרוא את המחשבות שלנו? האם זה עדיין בגדר החלום או שאנחנו כבר מתקרבים לפתרון? היום ב-shorthebrewpapereviews# סוקרים מאמר שבנה מודל לחיזוי (שחזור) תמונה שמראים לאדם מאות electroencephalogram (EEG) המוקלט מהמוח שלו. "

Set-RangeText "היום ב-#shorthebrewpapereviews סוקרים מאמר המציע גישה לזיהוי פרומפטים זדוניים שעלולים לגרום ליצירת תוכן מסוכן. המאמר מציע לבנות סטים של פרומפטים זדוניים הממקסמים 3 מטריקות שכל אחת מהם מודדת היבט שונה של ״זדוניות״ הפרומפטים מהסט הזה. היעד הראשון הוא מקסום סבירות של יצירת תוכן מסוכן עם פרומפטים מהסט, השני הוא הגיוון הסמנטי של הפרומפטים (כמה שפחות דמיון בין הפרומפטים) והיעד השלישי הוא הנראות ״הטובה״ של פרומפטים אלו (כלומר העדר של מילים גסות או בעלות תוכן מיני מובהק)." "המאמר מאמץ גישה משולבת לעיבוד של אות EEG: מצד אחד מנסים להפיק מהאות פיצ'רים עדינים(fine-grained) של התמונה בדמות מפת בולטות (saliency map) המפיקה את הפיצ'רים הויזואלים החשובים של התמונה (silhouette). "

Set-RangeText "המאמר משתמש במודל שפה בשביל ליצור פרומפטים אלו באמצעות מנגנון למידה in-context. האלגוריתם מתחיל בכמה פרומפטים זדוניים שנכתבו על ידי בני אדם ואז משתמשים במודל שפה כדי לגנרט פרומפטים זדוניים באמצעות מודל שפה (למידה in-context). עבור כל פרומפט זדוני שהצליח (יצר תוכן מסוכן) יוצרים סטים שבהם כל פרומפט מהסט מוחלף בפרומפט החדש ובוחרים מהם את הסט שממקסם לנו פונקציית היעד." "מצד שני מפיקים מהאות גם את הפיצ'רים הגסים של התמונה (ייצוג הכותרת שלה). שני הפיצ'רים אלו מזינים למודל דיפוזיה לטנטי (כמו Stable Diffusion) שמטרתו לשחזר את התמונה. הפיצ'רים העדינים (מפת בולטות) מחושבת בשני שלבים. "

# Last <w:t> is replaced AND split into two <w:t> runs (a brand new
# <w:br/><w:br/> pair is introduced) -- embed $lb twice in the replacement.
$newSplit = "בשלב הראשון מחשבים את הייצוג הלטנטי של אות ה-EEG עם למידה ניגודית (מקרבים ייצוגים של אותות EEG לתמונות דומות ומקרבים את אלו לתמונות לא דומות). בשלב השני מאמנים GAN מבוסס על hinge loss (כן עדיין משתמשים בהם) כדי ליצור מפת בולטות של התמונה (הדגימות ה״אמיתיות״ כאן הן התמונות שמראים אותן לאנשים). " + $lb + $lb + "הפיצ'רים הגסים מחושבים באופן הבא: יוצרים כותרת של התמונה עם מודל מאומן BLIP (מוקפאת) ומעבירים דרך CLIP כדי ליצור את ייצוגה. ואז מאמנים מודל כך הייצוג הגס המופק מהאות יהיה קרוב לייצוג של כותרת התמונה. ואז מכניסים את מפת הבולטות יחד עם ייצוג הכותרת של התמונה למודל דיפוזיה לטנטי כדי לשחזר את התמונה (האנקודר והדקודר מוקפאים). זה כל הקסם בגדול…"
Set-RangeText "איך יודעים שפרומפט הצליח ליצור תוכן לא ראוי? משתמשים במודלים מאומנים לזיהוי תוכן לא בטוח (כמו NudeNet או Q16)." $newSplit

# --- Trailing empty paragraph: Heading2 -> Normal style ---------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Style.NameLocal -eq "Heading 2") {
        $p.Style = "Normal"
    }
}

Write-Output "done"
